$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact string representation (avoid Excel
# auto-converting numeric-looking strings like '2.30' -> 2.3 or '1.00' -> 1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.934.99'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.248.81'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.34'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.22'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.90%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +4.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0982'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.35'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.47'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.105'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.583.71'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.51'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.08'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.831'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.245.37'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.47%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '43.861.21'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0985'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +4.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.08'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.49'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.85%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.45%  '
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Toncoin'
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.33'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +21.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.83'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.57%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '173.37'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.06%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.42'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.78%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.58%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.93'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.66'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.40%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.39'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.28'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.14%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.34%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.61'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.57%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.88'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.03'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0945'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.87%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.99%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.442.40'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.78%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.76%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.83'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -11.12%  '
